$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 242.3
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H40").Value = 4352.533
$ws.Range("I40").Value = 2915.3333
$ws.Range("K40").Value = 2915.3333
$ws.Range("M40").Value = -2740.3333
$ws.Range("H41").Value = 403.33334
$ws.Range("I41").Value = 522.625
$ws.Range("J41").Value = 164.75
$ws.Range("K41").Value = 522.625
$ws.Range("L41").Value = 164.75
$ws.Range("M41").Value = -82.625
$ws.Range("N41").Value = -1044.75
$ws.Range("H64").Value = 5340.2
$ws.Range("J64").Value = 5340.2
$ws.Range("L64").Value = 5340.2
$ws.Range("N64").Value = -5836.2
$ws.Range("H67").Value = 5340.2
$ws.Range("J67").Value = 5340.2
$ws.Range("L67").Value = 5340.2
$ws.Range("N67").Value = -7056.2
$ws.Range("H70").Value = 2661.111
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 2661.111
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7983.333
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8523.332999999999
$ws.Range("H73").Value = 2661.111
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2661.111
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7983.333
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9855.332999999999
$ws.Range("H74").Value = 5398.8
$ws.Range("I74").Value = 4887.5557
$ws.Range("K74").Value = 4887.5557
$ws.Range("M74").Value = -3951.5557
$ws.Range("H77").Value = 5398.8
$ws.Range("I77").Value = 4887.5557
$ws.Range("K77").Value = 24437.7785
$ws.Range("M77").Value = -19757.7785
$ws.Range("H80").Value = 29418.285
$ws.Range("I80").Value = 75567.875
$ws.Range("J80").Value = 1018.53845
$ws.Range("K80").Value = 226703.625
$ws.Range("L80").Value = 3055.61535
$ws.Range("M80").Value = -225705.625
$ws.Range("N80").Value = -5051.61535
$ws.Range("H83").Value = 29418.285
$ws.Range("I83").Value = 75567.875
$ws.Range("J83").Value = 1018.53845
$ws.Range("K83").Value = 680110.875
$ws.Range("L83").Value = 9166.84605
$ws.Range("M83").Value = -675118.875
$ws.Range("N83").Value = -19150.84605
$ws.Range("H96").Value = 4683.3335
$ws.Range("I96").Value = 900
$ws.Range("K96").Value = 2700
$ws.Range("M96").Value = -1327
$ws.Range("H98").Value = 1389.3846
$ws.Range("I98").Value = 1389.3846
$ws.Range("K98").Value = 1389.3846
$ws.Range("M98").Value = 108.6153999999999
$ws.Range("H113").Value = 4371.5713
$ws.Range("I113").Value = 3519.125
$ws.Range("J113").Value = 5508.1665
$ws.Range("K113").Value = 3519.125
$ws.Range("L113").Value = 5508.1665
$ws.Range("M113").Value = -265.125
$ws.Range("N113").Value = -12016.1665
$ws.Range("H115").Value = 1249.9
$ws.Range("I115").Value = 269.5
$ws.Range("K115").Value = 808.5
$ws.Range("M115").Value = 758.5
$ws.Range("H122").Value = 1389.3846
$ws.Range("I122").Value = 1389.3846
$ws.Range("K122").Value = 4168.1538
$ws.Range("M122").Value = -1718.1538
$ws.Range("H132").Value = 3868.077
$ws.Range("I132").Value = 3174.111
$ws.Range("J132").Value = 12195.667
$ws.Range("K132").Value = 9522.332999999999
$ws.Range("L132").Value = 36587.001
$ws.Range("M132").Value = -6992.332999999999
$ws.Range("N132").Value = -41647.001
$ws.Range("H137").Value = 5745.75
$ws.Range("I137").Value = 1229.2222
$ws.Range("J137").Value = 19295.334
$ws.Range("K137").Value = 3687.6666
$ws.Range("L137").Value = 57886.00199999999
$ws.Range("M137").Value = -1137.6666
$ws.Range("N137").Value = -62986.00199999999
$ws.Range("H138").Value = 5341.6787
$ws.Range("I138").Value = 4717.7
$ws.Range("J138").Value = 5477.326
$ws.Range("K138").Value = 14153.1
$ws.Range("L138").Value = 16431.978
$ws.Range("M138").Value = -9013.099999999999
$ws.Range("N138").Value = -26711.978
$ws.Range("H141").Value = 5967.731
$ws.Range("I141").Value = 6226.5
$ws.Range("J141").Value = 2862.5
$ws.Range("K141").Value = 18679.5
$ws.Range("L141").Value = 8587.5
$ws.Range("M141").Value = -13499.5
$ws.Range("N141").Value = -18947.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 9004.5
$ws.Range("I26").Value = 9004.5
$ws.Range("K26").Value = 9004.5
$ws.Range("M26").Value = -8674.5
$ws.Range("H32").Value = 25849.54
$ws.Range("I32").Value = 27895.162
$ws.Range("K32").Value = 27895.162
$ws.Range("M32").Value = -27608.162
$ws.Range("H45").Value = 2468.889
$ws.Range("I45").Value = 911.9091
$ws.Range("K45").Value = 911.9091
$ws.Range("M45").Value = -534.9091
$ws.Range("H61").Value = 3615
$ws.Range("I61").Value = 3371.3333
$ws.Range("J61").Value = 4199.8
$ws.Range("K61").Value = 3371.3333
$ws.Range("L61").Value = 4199.8
$ws.Range("M61").Value = -3159.3333
$ws.Range("N61").Value = -4623.8
$ws.Range("H74").Value = 173336.67
$ws.Range("I74").Value = 173336.67
$ws.Range("K74").Value = 173336.67
$ws.Range("M74").Value = -172462.67
$ws.Range("H77").Value = 173336.67
$ws.Range("I77").Value = 173336.67
$ws.Range("K77").Value = 866683.3500000001
$ws.Range("M77").Value = -862315.3500000001
$ws.Range("H110").Value = 2997.037
$ws.Range("I110").Value = 3207.2104
$ws.Range("J110").Value = 2497.875
$ws.Range("K110").Value = 3207.2104
$ws.Range("L110").Value = 2497.875
$ws.Range("M110").Value = -1162.2104
$ws.Range("N110").Value = -6587.875
$ws.Range("H112").Value = 76666.664
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H122").Value = 1436.1177
$ws.Range("I122").Value = 1027.6
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 3082.8
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -632.7999999999997
$ws.Range("N122").Value = -18400
$ws.Range("H128").Value = 69999
$ws.Range("J128").Value = 69999
$ws.Range("L128").Value = 69999
$ws.Range("N128").Value = -79959
$ws.Range("H132").Value = 33537
$ws.Range("I132").Value = 43195.707
$ws.Range("K132").Value = 129587.121
$ws.Range("M132").Value = -127057.121
$ws.Range("H136").Value = 3615
$ws.Range("I136").Value = 3371.3333
$ws.Range("J136").Value = 4199.8
$ws.Range("K136").Value = 10113.9999
$ws.Range("L136").Value = 12599.4
$ws.Range("M136").Value = -7563.999899999999
$ws.Range("N136").Value = -17699.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 31350
$ws.Range("J6").Value = 31350
$ws.Range("L6").Value = 31350
$ws.Range("N6").Value = -31576
$ws.Range("H21").Value = 20013.666
$ws.Range("J21").Value = 20013.666
$ws.Range("L21").Value = 20013.666
$ws.Range("N21").Value = -20485.666
$ws.Range("H86").Value = 2437.375
$ws.Range("I86").Value = 2166.3333
$ws.Range("K86").Value = 2166.3333
$ws.Range("M86").Value = -1043.3333
$ws.Range("H89").Value = 2437.375
$ws.Range("I89").Value = 2166.3333
$ws.Range("K89").Value = 10831.6665
$ws.Range("M89").Value = -5215.666499999999
$ws.Range("H94").Value = 7312.316
$ws.Range("I94").Value = 8355.25
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 8355.25
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -7904.25
$ws.Range("N94").Value = -2652
$ws.Range("H99").Value = 63553.117
$ws.Range("I99").Value = 95036.63
$ws.Range("K99").Value = 95036.63
$ws.Range("M99").Value = -93538.63
$ws.Range("H105").Value = 4463.1177
$ws.Range("I105").Value = 4325
$ws.Range("J105").Value = 5499
$ws.Range("K105").Value = 4325
$ws.Range("L105").Value = 5499
$ws.Range("M105").Value = -2578
$ws.Range("N105").Value = -8993
$ws.Range("H107").Value = 2779.625
$ws.Range("I107").Value = 2539.6667
$ws.Range("J107").Value = 3499.5
$ws.Range("K107").Value = 2539.6667
$ws.Range("L107").Value = 3499.5
$ws.Range("M107").Value = -619.6667000000002
$ws.Range("N107").Value = -7339.5
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H134").Value = 1701.1875
$ws.Range("I134").Value = 1701.1875
$ws.Range("K134").Value = 5103.5625
$ws.Range("M134").Value = -2568.5625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5447.7393
$ws.Range("I31").Value = 3251.1
$ws.Range("J31").Value = 7137.4614
$ws.Range("K31").Value = 3251.1
$ws.Range("L31").Value = 7137.4614
$ws.Range("M31").Value = -2956.1
$ws.Range("N31").Value = -7727.4614
$ws.Range("H34").Value = 5447.7393
$ws.Range("I34").Value = 3251.1
$ws.Range("J34").Value = 7137.4614
$ws.Range("K34").Value = 3251.1
$ws.Range("L34").Value = 7137.4614
$ws.Range("M34").Value = -3049.1
$ws.Range("N34").Value = -7541.4614
$ws.Range("H62").Value = 2874.5
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 2874.5
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H86").Value = 17159.58
$ws.Range("I86").Value = 32555.072
$ws.Range("J86").Value = 4480.9414
$ws.Range("K86").Value = 32555.072
$ws.Range("L86").Value = 4480.9414
$ws.Range("M86").Value = -31432.072
$ws.Range("N86").Value = -6726.9414
$ws.Range("H89").Value = 17159.58
$ws.Range("I89").Value = 32555.072
$ws.Range("J89").Value = 4480.9414
$ws.Range("K89").Value = 162775.36
$ws.Range("L89").Value = 22404.707
$ws.Range("M89").Value = -157159.36
$ws.Range("N89").Value = -33636.70699999999
$ws.Range("H107").Value = 503.5
$ws.Range("I107").Value = 458.36365
$ws.Range("K107").Value = 458.36365
$ws.Range("M107").Value = 1461.63635

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1947.3334
$ws.Range("I7").Value = 143
$ws.Range("K7").Value = 429
$ws.Range("M7").Value = -317
$ws.Range("H92").Value = 585.9167
$ws.Range("J92").Value = 859.25
$ws.Range("L92").Value = 2577.75
$ws.Range("N92").Value = -5073.75
$ws.Range("H105").Value = 9166.333000000001
$ws.Range("J105").Value = 9166.333000000001
$ws.Range("L105").Value = 27498.999
$ws.Range("N105").Value = -32740.999
$ws.Range("H113").Value = 770.30304
$ws.Range("J113").Value = 1029.7222
$ws.Range("L113").Value = 3089.1666
$ws.Range("N113").Value = -7429.1666
$ws.Range("H122").Value = 372.42105
$ws.Range("I122").Value = 339.27274
$ws.Range("J122").Value = 418
$ws.Range("K122").Value = 3053.45466
$ws.Range("L122").Value = 3762
$ws.Range("M122").Value = -603.4546599999999
$ws.Range("N122").Value = -8662
$ws.Range("H131").Value = 2133231.2
$ws.Range("J131").Value = 2864002.5
$ws.Range("L131").Value = 8592007.5
$ws.Range("N131").Value = -8602087.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2500
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2884
$ws.Range("N3").Value = -2232
$ws.Range("H10").Value = 20634
$ws.Range("J10").Value = 29999
$ws.Range("L10").Value = 29999
$ws.Range("N10").Value = -30337
$ws.Range("H87").Value = 69999
$ws.Range("J87").Value = 69999
$ws.Range("L87").Value = 69999
$ws.Range("N87").Value = -72495
$ws.Range("H90").Value = 69999
$ws.Range("J90").Value = 69999
$ws.Range("L90").Value = 209997
$ws.Range("N90").Value = -222477
$ws.Range("H102").Value = 2743.842
$ws.Range("I102").Value = 2147.4375
$ws.Range("J102").Value = 5924.6665
$ws.Range("K102").Value = 2147.4375
$ws.Range("L102").Value = 5924.6665
$ws.Range("M102").Value = -525.4375
$ws.Range("N102").Value = -9168.666499999999
$ws.Range("H107").Value = 57046.777
$ws.Range("I107").Value = 84491.086
$ws.Range("K107").Value = 84491.086
$ws.Range("M107").Value = -82571.086
$ws.Range("H122").Value = 3206.6316
$ws.Range("I122").Value = 3211.7334
$ws.Range("J122").Value = 3187.5
$ws.Range("K122").Value = 9635.200199999999
$ws.Range("L122").Value = 9562.5
$ws.Range("M122").Value = -7185.200199999999
$ws.Range("N122").Value = -14462.5
$ws.Range("H126").Value = 5317.55
$ws.Range("I126").Value = 4471.7144
$ws.Range("J126").Value = 7291.1665
$ws.Range("K126").Value = 13415.1432
$ws.Range("L126").Value = 21873.4995
$ws.Range("M126").Value = -10945.1432
$ws.Range("N126").Value = -26813.4995
$ws.Range("H131").Value = 49984
$ws.Range("J131").Value = 49984
$ws.Range("L131").Value = 49984
$ws.Range("N131").Value = -60064
$ws.Range("H132").Value = 62777.41
$ws.Range("I132").Value = 85956.336
$ws.Range("K132").Value = 257869.008
$ws.Range("M132").Value = -255339.008

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2983
$ws.Range("I40").Value = 2330.875
$ws.Range("K40").Value = 2330.875
$ws.Range("M40").Value = -2194.875
$ws.Range("H68").Value = 5449.6665
$ws.Range("I68").Value = 4883.3335
$ws.Range("J68").Value = 5732.8335
$ws.Range("K68").Value = 4883.3335
$ws.Range("L68").Value = 5732.8335
$ws.Range("M68").Value = -4134.3335
$ws.Range("N68").Value = -7230.8335
$ws.Range("H71").Value = 5449.6665
$ws.Range("I71").Value = 4883.3335
$ws.Range("J71").Value = 5732.8335
$ws.Range("K71").Value = 24416.6675
$ws.Range("L71").Value = 28664.1675
$ws.Range("M71").Value = -20672.6675
$ws.Range("N71").Value = -36152.1675
$ws.Range("H82").Value = 2294.5625
$ws.Range("I82").Value = 809
$ws.Range("J82").Value = 3072.7144
$ws.Range("K82").Value = 809
$ws.Range("L82").Value = 3072.7144
$ws.Range("M82").Value = -448
$ws.Range("N82").Value = -3794.7144
$ws.Range("H85").Value = 2294.5625
$ws.Range("I85").Value = 809
$ws.Range("J85").Value = 3072.7144
$ws.Range("K85").Value = 809
$ws.Range("L85").Value = 3072.7144
$ws.Range("M85").Value = 439
$ws.Range("N85").Value = -5568.7144
$ws.Range("H93").Value = 1611.3158
$ws.Range("J93").Value = 3056.8572
$ws.Range("L93").Value = 3056.8572
$ws.Range("N93").Value = -5552.8572
$ws.Range("H100").Value = 3538.9546
$ws.Range("I100").Value = 3299
$ws.Range("K100").Value = 3299
$ws.Range("M100").Value = -2758
$ws.Range("H125").Value = 74999.5
$ws.Range("J125").Value = 74999.5
$ws.Range("L125").Value = 74999.5
$ws.Range("N125").Value = -84839.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2004
$ws.Range("I7").Value = 2004
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2004
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1891
$ws.Range("N7").ClearContents()
$ws.Range("H62").Value = 93039.914
$ws.Range("J62").Value = 156482.58
$ws.Range("L62").Value = 156482.58
$ws.Range("N62").Value = -157730.58
$ws.Range("H65").Value = 93039.914
$ws.Range("J65").Value = 156482.58
$ws.Range("L65").Value = 782412.8999999999
$ws.Range("N65").Value = -788652.8999999999
$ws.Range("H107").Value = 622.63635
$ws.Range("I107").Value = 622.63635
$ws.Range("K107").Value = 1867.90905
$ws.Range("M107").Value = 52.09095000000002
$ws.Range("H122").Value = 320.6111
$ws.Range("I122").Value = 450.88235
$ws.Range("K122").Value = 1352.64705
$ws.Range("M122").Value = 1097.35295
$ws.Range("H131").Value = 43978.668
$ws.Range("J131").Value = 43978.668
$ws.Range("L131").Value = 43978.668
$ws.Range("N131").Value = -54058.668
$ws.Range("H132").Value = 60095
$ws.Range("I132").Value = 66117.78
$ws.Range("K132").Value = 198353.34
$ws.Range("M132").Value = -195823.34
$ws.Range("H136").Value = 4559.787
$ws.Range("I136").Value = 4680.846
$ws.Range("K136").Value = 14042.538
$ws.Range("M136").Value = -11492.538
$ws.Range("H140").Value = 78724.5
$ws.Range("J140").Value = 78724.5
$ws.Range("L140").Value = 78724.5
$ws.Range("N140").Value = -89084.5
$ws.Range("H141").Value = 86499.25
$ws.Range("J141").Value = 86499.25
$ws.Range("L141").Value = 86499.25
$ws.Range("N141").Value = -96859.25
